# celebrity-zoo.docx edit
#   1. Insert a new "Meta description" paragraph right after the H1 title.
#   2. Near the end of the doc, drop the duplicated bold title paragraph and
#      replace the italic "meta description" paragraph's text with the new
#      DALL-E image prompt (keeping its leading empty run + italic run).

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParaText($para) {
    # Paragraph.Range.Text includes the trailing paragraph-mark char(s);
    # strip them so we can compare against plain text.
    return $para.Range.Text.TrimEnd([char]13, [char]7)
}

$titleText = "Play Celebrity Zoo Free - Review of Unique and Playful Slot Game"

# ---------------------------------------------------------------------
# Step 1: right after the FIRST paragraph carrying the title (the H1
# heading), insert the new meta-description paragraph:
#   <empty run><bold "Meta description"><normal rest-of-text>
# ---------------------------------------------------------------------
$titlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ((Get-ParaText $d.Paragraphs($i)) -eq $titleText) {
        $titlePara = $d.Paragraphs($i)
        break
    }
}

$titlePara.Range.InsertParagraphAfter() | Out-Null
$insertPoint = $d.Paragraphs($titlePara.Index + 1).Range
$insertPoint.Collapse(1)
$metaXml = "<w:p $wNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Celebrity Zoo is an entertaining slot game featuring celebrity animals, free spins, and two bonus games. Try it for free today!</w:t></w:r></w:p>"
$insertPoint.InsertXML($metaXml)

# ---------------------------------------------------------------------
# Step 2: near the end, the title text is duplicated in its own (bold)
# paragraph, immediately followed by an italic "meta description"
# paragraph. Find the LAST paragraph with the title text (the duplicate,
# not the H1 at the top) and delete it, then rewrite the text of the
# italic paragraph that follows it with the new DALL-E prompt.
# ---------------------------------------------------------------------
$boldParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ((Get-ParaText $d.Paragraphs($i)) -eq $titleText) {
        $boldParaIndex = $i
    }
}

$italicParaIndex = $boldParaIndex + 1
$d.Paragraphs($boldParaIndex).Range.Delete()
# after deleting the bold paragraph, every later paragraph shifts down by one
$italicPara = $d.Paragraphs($italicParaIndex - 1)

$promptText = 'Prompt: "Create a cartoon-style image for the game ''Celebrity Zoo'' that features a happy Maya warrior with glasses." DALLE, I need you to create a feature image for the slot game ''Celebrity Zoo'' that highlights its quirky and playful atmosphere. The image should be in a cartoon style that features a happy Maya warrior with glasses. The Maya warrior should be depicted in bright colors to match the game''s whimsical graphics, with a big smile on his face to convey the game''s fun and lighthearted tone. The Maya warrior should also be holding a camera to symbolize the game''s Paparazzi Bonus Game. This feature image should be eye-catching and vibrant to catch the players'' attention and make them want to try out the game. Make sure the image is designed to fit the game''s theme and can convey the game''s exciting features.'
$promptXml = "<w:p $wNs><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>$promptText</w:t></w:r></w:p>"
$italicPara.Range.InsertXML($promptXml)

Write-Output "Edit complete"
